# "added ifoCAST full series evaluation"
#
# The sheet holds a staircase-shaped matrix of QoQ forecast errors: each
# data row (vintages 10..28, one per row r=2..20) holds a run of numeric
# cells starting at column B. Evaluating one more ifoCAST data point per
# vintage rolls the whole window: the oldest (leftmost, column B) error
# for every row is dropped and every remaining value shifts one column
# to the left. Rows that were already at the sheet's full width (B:K)
# get a freshly evaluated error appended at the new rightmost column K;
# rows that were already narrower than B:K simply shrink by one column,
# since no replacement value is available yet for those shorter
# vintages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly evaluated right-hand (column K) errors for the full-width rows
# (r = 2..10), keyed by row number.
$newK = @{
    2  = 0.2199829514341669
    3  = -1.427026823174395
    4  = 1.345095091002794
    5  = 1.835270244654998
    6  = -0.1485141439230462
    7  = 0.5580467639488803
    8  = -0.4310464000952693
    9  = 0.5995033638472159
    10 = -0.5018523531907899
}

# Rows 2..10 are full width (columns B..K = 2..11): shift every value one
# column to the left, then drop in the newly evaluated column K value.
# NOTE: `.Value` must be invoked with `()` on the read side - without it,
# the COM getter isn't actually called.
foreach ($r in 2..10) {
    for ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($r, $c + 1).Value()
    }
    $ws.Cells.Item($r, 11).Value = $newK[$r]
}

# Rows 11..20 were already narrower than B:K. Each row's last used column
# (1-based) before this edit:
$rowEnd = @{
    11 = 11   # B11:K11 -> B11:J11
    12 = 10   # B12:J12 -> B12:I12
    13 = 9    # B13:I13 -> B13:H13
    14 = 8    # B14:H14 -> B14:G14
    15 = 7    # B15:G15 -> B15:F15
    16 = 6    # B16:F16 -> B16:E16
    17 = 5    # B17:E17 -> B17:D17
    18 = 4    # B18:D18 -> B18:C18
    19 = 3    # B19:C19 -> B19:B19
    20 = 2    # B20      -> (empty)
}

foreach ($r in 11..20) {
    $endCol = $rowEnd[$r]
    for ($c = 2; $c -lt $endCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($r, $c + 1).Value()
    }
    $ws.Cells.Item($r, $endCol).ClearContents()
}
